$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.133.35"
$ws.Range("E2").Value = "  -0.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.270.83"
$ws.Range("E3").Value = "  +0.37%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.74"
$ws.Range("E5").Value = "  +0.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.42"
$ws.Range("E6").Value = "  -0.45%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.130"
$ws.Range("E9").Value = "  -2.16%  "

# Row 10
$ws.Range("E10").Value = "  -0.46%  "

# Row 11
$ws.Range("E11").Value = "  -2.93%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.838.66"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13
$ws.Range("E13").Value = "  +1.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.126.40"
$ws.Range("E14").Value = "  -0.24%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.38"
$ws.Range("E15").Value = "  -3.09%  "

# Row 16
$ws.Range("E16").Value = "  -2.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.272.67"
$ws.Range("E17").Value = "  +0.93%  "

# Row 18
$ws.Range("E18").Value = "  -2.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.29"
$ws.Range("E19").Value = "  -2.49%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "416.70"
$ws.Range("E20").Value = "  +5.89%  "

# Row 21
$ws.Range("E21").Value = "  -2.23%  "

# Row 22
$ws.Range("E22").Value = "  +0.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.06"
$ws.Range("E23").Value = "  -0.40%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.508"
$ws.Range("E24").Value = "  -2.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000118"
$ws.Range("E25").Value = "  -2.17%  "

# Row 26
$ws.Range("E26").Value = "  -0.96%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.36"
$ws.Range("E27").Value = "  -4.64%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("E29").Value = "  -1.78%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.62"
$ws.Range("E30").Value = "  -1.45%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.42"
$ws.Range("E31").Value = "  -5.21%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.85"
$ws.Range("E32").Value = "  -4.50%  "

# Row 33
$ws.Range("E33").Value = "  -3.45%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "164.55"
$ws.Range("E34").Value = "  +1.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.44"
$ws.Range("E35").Value = "  -4.69%  "

# Row 36
$ws.Range("E36").Value = "  -4.34%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.71"
$ws.Range("E37").Value = "  -0.62%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.792"
$ws.Range("E38").Value = "  -3.87%  "

# Row 39
$ws.Range("E39").Value = "  -3.53%  "

# Row 40
$ws.Range("E40").Value = "  -4.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.630.98"
$ws.Range("E41").Value = "  -0.92%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0675"
$ws.Range("E42").Value = "  -2.30%  "

# Row 43
$ws.Range("E43").Value = "  -4.20%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "335.41"
$ws.Range("E44").Value = "  -1.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.24"
$ws.Range("E45").Value = "  -4.66%  "

# Row 46
$ws.Range("E46").Value = "  -3.18%  "

# Row 47
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.987"
$ws.Range("E47").Value = "  -0.43%  "

# Row 48
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.23"
$ws.Range("E48").Value = "  -2.16%  "

# Row 49
$ws.Range("E49").Value = "  -1.72%  "

# Row 50
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "30.65"
$ws.Range("E51").Value = "  -2.93%  "
